$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new D (Price) text, new E (Volume(1h)) text.
# $null means "leave this cell unchanged".
# Values that look like plain numbers are prefixed with a leading
# apostrophe so Excel keeps them as literal text (matching the original
# file's layout where Price/Volume are stored as text, e.g. "1.00" would
# otherwise collapse to the number 1).
$updates = @(
    @{ Row = 2;  D = "70.266.71";   E = "  +0.45%  " },
    @{ Row = 3;  D = "3.608.02";    E = "  +2.24%  " },
    @{ Row = 4;  D = $null;         E = "  +0.07%  " },
    @{ Row = 5;  D = "'604.40";     E = "  +0.05%  " },
    @{ Row = 6;  D = "'196.37";     E = "  -0.10%  " },
    @{ Row = 7;  D = $null;         E = "  +0.07%  " },
    @{ Row = 8;  D = "'1.00";       E = $null },
    @{ Row = 9;  D = $null;         E = "  -1.75%  " },
    @{ Row = 10; D = $null;         E = "  -1.20%  " },
    @{ Row = 11; D = "'53.88";      E = "  -0.41%  " },
    @{ Row = 12; D = $null;         E = "  +0.30%  " },
    @{ Row = 14; D = "4.180.08";    E = "  +2.30%  " },
    @{ Row = 15; D = "'13.12";      E = "  +3.53%  " },
    @{ Row = 16; D = "'599.12";     E = "  -0.54%  " },
    @{ Row = 17; D = "70.376.54";   E = "  +0.41%  " },
    @{ Row = 19; D = "3.608.52";    E = "  +2.06%  " },
    @{ Row = 20; D = $null;         E = "  +1.41%  " },
    @{ Row = 21; D = $null;         E = "  +0.01%  " },
    @{ Row = 22; D = "'17.82";      E = "  -2.29%  " },
    @{ Row = 23; D = "'5.18";       E = "  -2.12%  " },
    @{ Row = 24; D = "'102.12";     E = "  -1.35%  " },
    @{ Row = 25; D = $null;         E = "  +0.02%  " },
    @{ Row = 26; D = "'3.03";       E = "  -2.23%  " },
    @{ Row = 27; D = "'10.76";      E = "  -1.85%  " },
    @{ Row = 28; D = $null;         E = "  -0.73%  " },
    @{ Row = 29; D = "'33.85";      E = "  +0.57%  " },
    @{ Row = 30; D = "'4.77";       E = "  +5.54%  " },
    @{ Row = 31; D = "'7.17";       E = "  +0.60%  " },
    @{ Row = 32; D = "'12.29";      E = "  -3.54%  " },
    @{ Row = 33; D = $null;         E = "  +0.67%  " },
    @{ Row = 34; D = $null;         E = "  -0.35%  " },
    @{ Row = 35; D = "0.0₃0896";    E = "  +7.67%  " },
    @{ Row = 36; D = "3.909.04";    E = "  +4.55%  " },
    @{ Row = 37; D = "'3.10";       E = "  +0.24%  " },
    @{ Row = 38; D = $null;         E = "  +0.10%  " },
    @{ Row = 39; D = "'520.71";     E = "  +5.73%  " },
    @{ Row = 40; D = "'36.96";      E = "  +0.30%  " },
    @{ Row = 41; D = $null;         E = "  -1.14%  " },
    @{ Row = 42; D = $null;         E = "  -2.08%  " },
    @{ Row = 43; D = $null;         E = "  -2.13%  " },
    @{ Row = 44; D = "'0.0454";     E = "  -0.70%  " },
    @{ Row = 45; D = "'3.42";       E = "  +2.85%  " },
    @{ Row = 46; D = "'2.86";       E = "  +0.77%  " },
    @{ Row = 47; D = $null;         E = "  -0.04%  " },
    @{ Row = 48; D = $null;         E = "  -0.43%  " },
    @{ Row = 49; D = $null;         E = "  -0.29%  " },
    @{ Row = 50; D = "'0.000252";   E = "  +2.43%  " },
    @{ Row = 51; D = $null;         E = "  +0.07%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $ws.Cells.Item($r, 4).Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}
